$wb = $excel.ActiveWorkbook

# Rename worksheets
$wb.Worksheets.Item("Include ValueSets").Name = "Include ValueSet #0"
$wb.Worksheets.Item("Include from CareSocialCodes").Name = "Include #1"

# Update metadata values on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.2.0"
$meta.Range("B8").Value = "2024-10-31T15:40:44+01:00"
